$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting the existing data (rows 2-31) down to (3-32)
$ws.Rows("2:2").Insert()

# Populate the newly inserted row 2 with the new insider-trade record
$ws.Range("A2").Value = "REFR"
$ws.Range("B2").Value = "Kaganowicz Alexander"
$ws.Range("C2").Value = "Director"
$ws.Range("D2").Value = "Jun 14"
$ws.Range("E2").Value = "Buy"
$ws.Range("F2").Value = 1.94
$ws.Range("G2").Value = 2000
$ws.Range("H2").Value = 3880
$ws.Range("I2").Value = 164923
$ws.Range("J2").Value = "Jun 17 06:30 AM"
